# Sync attendance_reports: update "Recorded By" (column G) cell values
# so the recorder email lists are reordered to match the latest export.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$updates = @{
    "G2"  = "gehanadel@med.asu.edu.eg, System, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
    "G3"  = "eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, System, majorelle.magdy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
    "G4"  = "eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, gehanadel@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
    "G5"  = "eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
    "G6"  = "alshimaa.atef@med.asu.edu.egm, manar.montaser@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg"
    "G7"  = "Amera.a.saad@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg"
    "G12" = "Marina.youhana@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, dina.adel@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg"
    "G17" = "mohamed.saleem@med.asu.edu.eg, esraa.sami@med.asu.edu.eg"
    "G24" = "youstina.gamil@med.asu.edu.eg, Sarah.Mahdy@med.asu.edu.eg"
    "G25" = "Noran.Mahmoud@med.asu.edu.eg, menna-allah.gamil@med.asu.edu.eg"
    "G27" = "hana.amr@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"
    "G30" = "yassmen.ahmed@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
